$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- Title: "Challenge" -> "Challenges" ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Challenges"

# --- Content placeholder ---
$tf = $s.Shapes.Item(2).TextFrame
$tr = $tf.TextRange

$oldProblemsLink = "http://com1243.eecs.utk.edu:8888/crypto_warmup/"

# Remove the old "http://com1243.eecs.utk.edu:8888/crypto_warmup/" paragraph
# entirely (including its paragraph break); this shifts the "Solutions: ..."
# paragraph up to become paragraph 1. Its leading run ("Solutions: ") carries
# no hyperlink, which we use below as a clean insertion point so the new
# "Problems: " text we type does not inherit a hyperlink.
$tr.Paragraphs(1).Delete()

# Insert the new first line (plus a paragraph break) right before the clean
# "Solutions: " run so the new text stays unlinked.
$tr.Paragraphs(1).InsertBefore("Problems: " + [char]13)

$probPara = $tr.Paragraphs(1)
$probPara.InsertAfter("https://tiny.utk.edu/crypto-practice")

# "https://"  (8 chars) starts right after "Problems: " (10 chars)
$httpsPart = $probPara.Characters(11, 8)
$httpsPart.ActionSettings(1).Hyperlink.Address = $oldProblemsLink

# "tiny.utk.edu/crypto-practice" (29 chars) follows immediately after
$restPart = $probPara.Characters(19, 29)
$restPart.ActionSettings(1).Hyperlink.Address = $oldProblemsLink

# --- Solutions paragraph ---
$solPara = $tr.Paragraphs(2)

# Split "Solutions: " into "Solutions" + ": " (two clean runs)
$solWord = $solPara.Characters(1, 9)
$solWord.Text = "Solutions"
$colonSpace = $solPara.Characters(10, 2)
$colonSpace.Text = ": "

# Merge "https://" + "github.com/hackutk/historical-crypto" into a single run
$link = $solPara.Characters(12, 8 + 37)
$link.Text = "https://github.com/hackutk/historical-crypto"
